$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'20.427.23"
$ws.Range('E2').Value = '  +2.43%  '
$ws.Range('D3').Value = "'1.466.41"
$ws.Range('E3').Value = '  +4.45%  '
$ws.Range('E4').Value = '  +1.00%  '
$ws.Range('D5').Value = "'280.46"
$ws.Range('E5').Value = '  +2.69%  '
$ws.Range('D6').Value = "'0.8997"
$ws.Range('E6').Value = '  -10.22%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').Value = "'0.3183"
$ws.Range('E8').Value = '  +3.78%  '
$ws.Range('D9').Value = "'39.44"
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').Value = "'1.048"
$ws.Range('E10').Value = '  +5.52%  '
$ws.Range('D11').Value = "'0.06610"
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').Value = "'1.008"
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = "'5.544"
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').Value = "'17.87"
$ws.Range('E14').Value = '  +6.16%  '
$ws.Range('D15').Value = "'6.199"
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = "'1.472.90"
$ws.Range('E16').Value = '  +4.87%  '
$ws.Range('D17').Value = "'0.00001029"
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('D19').Value = "'0.8994"
$ws.Range('E19').Value = '  -10.31%  '
$ws.Range('D20').Value = "'70.48"
$ws.Range('E20').Value = '  -4.10%  '
$ws.Range('D21').Value = "'5.685"
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('D22').Value = "'14.58"
$ws.Range('E22').Value = '  +1.42%  '
$ws.Range('D23').Value = "'11.19"
$ws.Range('E23').Value = '  +3.78%  '
$ws.Range('D24').Value = "'2.287"
$ws.Range('E24').Value = '  -1.53%  '
$ws.Range('D25').Value = "'20.709.42"
$ws.Range('E25').Value = '  +3.89%  '
$ws.Range('D26').Value = "'2.266"
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('D27').Value = "'137.30"
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('D28').Value = "'17.42"
$ws.Range('E28').Value = '  +3.07%  '
$ws.Range('D29').Value = "'1.636.08"
$ws.Range('E29').Value = '  +4.44%  '
$ws.Range('D30').Value = "'113.20"
$ws.Range('E30').Value = '  +4.04%  '
$ws.Range('D31').Value = "'3.939"
$ws.Range('E31').Value = '  +2.35%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = "'0.8375"
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'5.111"
$ws.Range('E33').Value = '  -4.96%  '
$ws.Range('D34').Value = "'0.07800"
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('D35').Value = "'0.06088"
$ws.Range('E35').Value = '  +5.92%  '
$ws.Range('D36').Value = "'1.446"
$ws.Range('E36').Value = '  +13.95%  '
$ws.Range('D37').Value = "'4.844"
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').Value = "'1.158"
$ws.Range('E38').Value = '  +9.98%  '
$ws.Range('D39').Value = "'10.53"
$ws.Range('E39').Value = '  +3.17%  '
$ws.Range('D40').Value = "'0.02037"
$ws.Range('E40').Value = '  +0.54%  '
$ws.Range('D41').Value = "'0.1864"
$ws.Range('E41').Value = '  -2.55%  '
$ws.Range('D42').Value = "'0.9160"
$ws.Range('E42').Value = '  -8.49%  '
$ws.Range('D43').Value = "'0.5352"
$ws.Range('E43').Value = '  +1.49%  '
$ws.Range('D44').Value = "'3.579"
$ws.Range('E44').Value = '  +1.56%  '
$ws.Range('D45').Value = "'12.26"
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('D46').Value = "'6.844"
$ws.Range('E46').Value = '  -18.51%  '
$ws.Range('D47').Value = "'122.38"
$ws.Range('E47').Value = '  +12.06%  '
$ws.Range('D48').Value = "'0.5254"
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('D49').Value = "'1.820"
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('D50').Value = "'0.06424"
$ws.Range('E50').Value = '  +4.52%  '
$ws.Range('D51').Value = "'1.036"
$ws.Range('E51').Value = '  -0.97%  '
